# V. 117 "Blue Beetle"
# Adds a new movie entry "Blue Beetle" to the "Películas" sheet table,
# keeping the existing sort order (table is not re-sorted automatically)
# and moves the "latest addition" highlight style from row 36 (Damsel)
# to the newly inserted row.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Películas")

# ---------------------------------------------------------------------
# 1. Insert a block of blank rows right after row 36 (the row currently
#    holding the "latest addition" highlighted style). Because a plain
#    row Insert() duplicates formatting from the row immediately above,
#    every row in this new block inherits the highlighted style.
#    We insert exactly as many rows (96: rows 37-132) as needed so that,
#    after later removing the ones we don't want, the single surviving
#    highlighted row lands exactly on row 132 - the new row's final
#    target position.
# ---------------------------------------------------------------------
$ws.Range("37:132").Insert()

# ---------------------------------------------------------------------
# 2. The original rows 37-131 got pushed down to 133-227 by the insert
#    above. Copy their values/styles back into rows 37-131 so that,
#    content-wise, nothing really changed for them (row 132 is left
#    alone, still blank and still carrying the highlighted style).
# ---------------------------------------------------------------------
$ws.Range("B133:I227").Copy($ws.Range("B37:I131"))

# The bulk copy above only carries over literal values, so restore the
# per-row AVERAGE() formulas in column C for rows 37-131.
for ($r = 37; $r -le 131; $r++) {
    $ws.Cells.Item($r, 3).Formula = "=AVERAGE(D$r,E$r,E$r,F$r,G$r,H$r,H$r,I$r)"
}

# ---------------------------------------------------------------------
# 3. Rows 133-227 are now redundant duplicates of the original rows
#    37-131; delete them. This shifts the original rows 132-137 (which
#    are currently sitting at 228-233) back up to 133-138.
# ---------------------------------------------------------------------
$ws.Range("133:227").Delete()

# ---------------------------------------------------------------------
# 4. Row 132 is still blank, keeping the highlighted style inherited in
#    step 1. Fill it in with the "Blue Beetle" data.
# ---------------------------------------------------------------------
$ws.Cells.Item(132, 2).Value = "Blue Beetle"
$ws.Cells.Item(132, 4).Value = 3
$ws.Cells.Item(132, 5).Value = 3
$ws.Cells.Item(132, 6).Value = 2
$ws.Cells.Item(132, 7).Value = 4
$ws.Cells.Item(132, 8).Value = 5.9
$ws.Cells.Item(132, 9).Value = 4.8
$ws.Cells.Item(132, 3).Formula = "=AVERAGE(D132,E132,E132,F132,G132,H132,H132,I132)"

# ---------------------------------------------------------------------
# 5. Row 36 (Damsel) is no longer the latest addition; drop its
#    highlighted style back to the regular left-aligned style used by
#    every other row, copying the format from the row above it.
# ---------------------------------------------------------------------
$ws.Range("B35").Copy()
$ws.Range("B36").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# ---------------------------------------------------------------------
# 6. Resize the table / autofilter to include the new row.
# ---------------------------------------------------------------------
$lo = $ws.ListObjects.Item(1)
$lo.Resize($ws.Range("B2:I138"))

# ---------------------------------------------------------------------
# 7. Update the view so the active cell matches the new last row, and
#    keep the same scroll position as before.
# ---------------------------------------------------------------------
$ws.Range("C138").Select()
